$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 63
$ws_ALC.Range("H63").Value = 34666.332
$ws_ALC.Range("J63").Value = 34999.5
$ws_ALC.Range("L63").Value = 34999.5
$ws_ALC.Range("N63").Value = -36247.5

# ALC row 66
$ws_ALC.Range("H66").Value = 34666.332
$ws_ALC.Range("J66").Value = 34999.5
$ws_ALC.Range("L66").Value = 104998.5
$ws_ALC.Range("N66").Value = -111238.5

# ALC row 113
$ws_ALC.Range("H113").Value = 113131.664
$ws_ALC.Range("I113").Value = 144883.58
$ws_ALC.Range("J113").Value = 2000
$ws_ALC.Range("K113").Value = 144883.58
$ws_ALC.Range("L113").Value = 2000
$ws_ALC.Range("M113").Value = -141629.58
$ws_ALC.Range("N113").Value = -8508

# ALC row 129
$ws_ALC.Range("H129").Value = 2117
$ws_ALC.Range("J129").Value = 837.4693600000001
$ws_ALC.Range("L129").Value = 2512.40808
$ws_ALC.Range("N129").Value = -12512.40808

# ALC row 132
$ws_ALC.Range("H132").Value = 6762926
$ws_ALC.Range("I132").Value = 7819161
$ws_ALC.Range("J132").Value = 3020
$ws_ALC.Range("K132").Value = 23457483
$ws_ALC.Range("L132").Value = 9060
$ws_ALC.Range("M132").Value = -23454953
$ws_ALC.Range("N132").Value = -14120

# ALC row 138
$ws_ALC.Range("H138").Value = 2039.4242
$ws_ALC.Range("I138").Value = 3129.125
$ws_ALC.Range("J138").Value = 1690.72
$ws_ALC.Range("K138").Value = 9387.375
$ws_ALC.Range("L138").Value = 5072.16
$ws_ALC.Range("M138").Value = -4247.375
$ws_ALC.Range("N138").Value = -15352.16

# ARM row 32
$ws_ARM.Range("H32").Value = 24940.56
$ws_ARM.Range("I32").Value = 7065.0586
$ws_ARM.Range("K32").Value = 7065.0586
$ws_ARM.Range("M32").Value = -6778.0586

# ARM row 74
$ws_ARM.Range("H74").Value = 728.25
$ws_ARM.Range("I74").Value = 654.5599999999999
$ws_ARM.Range("J74").Value = 991.4286
$ws_ARM.Range("K74").Value = 654.5599999999999
$ws_ARM.Range("L74").Value = 991.4286
$ws_ARM.Range("M74").Value = 219.4400000000001
$ws_ARM.Range("N74").Value = -2739.4286

# ARM row 77
$ws_ARM.Range("H77").Value = 728.25
$ws_ARM.Range("I77").Value = 654.5599999999999
$ws_ARM.Range("J77").Value = 991.4286
$ws_ARM.Range("K77").Value = 3272.8
$ws_ARM.Range("L77").Value = 4957.143
$ws_ARM.Range("M77").Value = 1095.2
$ws_ARM.Range("N77").Value = -13693.143

# ARM row 122
$ws_ARM.Range("H122").Value = 1458.9333
$ws_ARM.Range("I122").Value = 1137.6
$ws_ARM.Range("K122").Value = 3412.8
$ws_ARM.Range("M122").Value = -962.7999999999997

# ARM row 132
$ws_ARM.Range("H132").Value = 4662.913
$ws_ARM.Range("I132").Value = 4825.1763
$ws_ARM.Range("K132").Value = 14475.5289
$ws_ARM.Range("M132").Value = -11945.5289

# BSM row 35
$ws_BSM.Range("H35").Value = 19483.25
$ws_BSM.Range("J35").Value = 19483.25
$ws_BSM.Range("L35").Value = 19483.25
$ws_BSM.Range("N35").Value = -20103.25

# BSM row 88
$ws_BSM.Range("H88").Value = 16666.5
$ws_BSM.Range("J88").Value = 16666.5
$ws_BSM.Range("L88").Value = 16666.5
$ws_BSM.Range("N88").Value = -17478.5

# BSM row 91
$ws_BSM.Range("H91").Value = 16666.5
$ws_BSM.Range("J91").Value = 16666.5
$ws_BSM.Range("L91").Value = 16666.5
$ws_BSM.Range("N91").Value = -19474.5

# CRP row 132
$ws_CRP.Range("H132").Value = 4152.3
$ws_CRP.Range("I132").Value = 4120.4707
$ws_CRP.Range("K132").Value = 12361.4121
$ws_CRP.Range("M132").Value = -9831.4121

# CUL row 55
$ws_CUL.Range("H55").Value = 15998.571
$ws_CUL.Range("J55").Value = 12168
$ws_CUL.Range("L55").Value = 36504
$ws_CUL.Range("N55").Value = -36858

# CUL row 131
$ws_CUL.Range("H131").Value = 1078.5143
$ws_CUL.Range("J131").Value = 1105.9104
$ws_CUL.Range("L131").Value = 3317.7312
$ws_CUL.Range("N131").Value = -13397.7312

# CUL row 132
$ws_CUL.Range("H132").Value = 1597.6857
$ws_CUL.Range("I132").Value = 742.0526
$ws_CUL.Range("J132").Value = 2613.75
$ws_CUL.Range("K132").Value = 6678.4734
$ws_CUL.Range("L132").Value = 23523.75
$ws_CUL.Range("M132").Value = -4148.4734
$ws_CUL.Range("N132").Value = -28583.75

# GSM row 43
$ws_GSM.Range("H43").Value = 5531.2856
$ws_GSM.Range("I43").Value = 1850
$ws_GSM.Range("J43").Value = 7003.8
$ws_GSM.Range("K43").Value = 1850
$ws_GSM.Range("L43").Value = 7003.8
$ws_GSM.Range("M43").Value = -1699
$ws_GSM.Range("N43").Value = -7305.8

# GSM row 46
$ws_GSM.Range("H46").Value = 11999.818
$ws_GSM.Range("I46").Value = 8000
$ws_GSM.Range("J46").Value = 12399.8
$ws_GSM.Range("K46").Value = 8000
$ws_GSM.Range("L46").Value = 12399.8
$ws_GSM.Range("M46").Value = -7844
$ws_GSM.Range("N46").Value = -12711.8

# GSM row 57
$ws_GSM.Range("H57").Value = 18933.334
$ws_GSM.Range("J57").Value = 18933.334
$ws_GSM.Range("L57").Value = 18933.334
$ws_GSM.Range("N57").Value = -20573.334

# GSM row 68
$ws_GSM.Range("H68").Value = 49000
$ws_GSM.Range("J68").Value = 49000
$ws_GSM.Range("L68").Value = 49000
$ws_GSM.Range("N68").Value = -50622

# GSM row 71
$ws_GSM.Range("H71").Value = 49000
$ws_GSM.Range("J71").Value = 49000
$ws_GSM.Range("L71").Value = 147000
$ws_GSM.Range("N71").Value = -155112

# GSM row 80
$ws_GSM.Range("H80").Value = 100105380
$ws_GSM.Range("I80").Value = 143005900
$ws_GSM.Range("J80").Value = 4133.3335
$ws_GSM.Range("K80").Value = 143005900
$ws_GSM.Range("L80").Value = 4133.3335
$ws_GSM.Range("M80").Value = -143004902
$ws_GSM.Range("N80").Value = -6129.3335

# GSM row 83
$ws_GSM.Range("H83").Value = 100105380
$ws_GSM.Range("I83").Value = 143005900
$ws_GSM.Range("J83").Value = 4133.3335
$ws_GSM.Range("K83").Value = 715029500
$ws_GSM.Range("L83").Value = 20666.6675
$ws_GSM.Range("M83").Value = -715024508
$ws_GSM.Range("N83").Value = -30650.6675

# GSM row 122
$ws_GSM.Range("H122").Value = 2840.4
$ws_GSM.Range("I122").Value = 2003
$ws_GSM.Range("J122").Value = 3398.6667
$ws_GSM.Range("K122").Value = 6009
$ws_GSM.Range("L122").Value = 10196.0001
$ws_GSM.Range("M122").Value = -3559
$ws_GSM.Range("N122").Value = -15096.0001

# GSM row 132
$ws_GSM.Range("H132").Value = 4372.273
$ws_GSM.Range("I132").Value = 1650
$ws_GSM.Range("K132").Value = 4950
$ws_GSM.Range("M132").Value = -2420

# LTW row 40
$ws_LTW.Range("H40").Value = 33533.125
$ws_LTW.Range("I40").Value = 47748.184
$ws_LTW.Range("K40").Value = 47748.184
$ws_LTW.Range("M40").Value = -47612.184

# LTW row 122
$ws_LTW.Range("H122").Value = 3200.818
$ws_LTW.Range("I122").Value = 3267.111
$ws_LTW.Range("J122").Value = 2902.5
$ws_LTW.Range("K122").Value = 9801.332999999999
$ws_LTW.Range("L122").Value = 8707.5
$ws_LTW.Range("M122").Value = -7351.332999999999
$ws_LTW.Range("N122").Value = -13607.5

# LTW row 132
$ws_LTW.Range("H132").Value = 5627.7646
$ws_LTW.Range("I132").Value = 6140.4
$ws_LTW.Range("K132").Value = 18421.2
$ws_LTW.Range("M132").Value = -15891.2

# WVR row 6
$ws_WVR.Range("H6").Value = 276951.5
$ws_WVR.Range("J6").Value = 276951.5
$ws_WVR.Range("L6").Value = 276951.5
$ws_WVR.Range("N6").Value = -277181.5

# WVR row 54
$ws_WVR.Range("H54").Value = 6924.7856
$ws_WVR.Range("J54").Value = 6913.615
$ws_WVR.Range("L54").Value = 6913.615
$ws_WVR.Range("N54").Value = -7953.615

# WVR row 81
$ws_WVR.Range("H81").Value = 333896.34
$ws_WVR.Range("I81").Value = 333866.34
$ws_WVR.Range("K81").Value = 667732.6800000001
$ws_WVR.Range("M81").Value = -666671.6800000001

# WVR row 84
$ws_WVR.Range("H84").Value = 333896.34
$ws_WVR.Range("I84").Value = 333866.34
$ws_WVR.Range("K84").Value = 3338663.4
$ws_WVR.Range("M84").Value = -3333359.4

# WVR row 132
$ws_WVR.Range("H132").Value = 2578.9302
$ws_WVR.Range("I132").Value = 2971.7144
$ws_WVR.Range("J132").Value = 1845.7333
$ws_WVR.Range("K132").Value = 8915.143199999999
$ws_WVR.Range("L132").Value = 5537.199900000001
$ws_WVR.Range("M132").Value = -6385.143199999999
$ws_WVR.Range("N132").Value = -10597.1999
